# Update Leve profit/price cells across sheets (scheduled runner sync)
$wb = $excel.ActiveWorkbook

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3016.9092
$ws.Range("J112").Value = 3424.625
$ws.Range("L112").Value = 10273.875
$ws.Range("N112").Value = -12489.875

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5402.4185
$ws.Range("I137").Value = 6886.846
$ws.Range("J137").Value = 3132.1177
$ws.Range("K137").Value = 20660.538
$ws.Range("L137").Value = 9396.3531
$ws.Range("M137").Value = -18110.538
$ws.Range("N137").Value = -14496.3531

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2904.5
$ws.Range("I138").Value = 1183.7142
$ws.Range("J138").Value = 3386.32
$ws.Range("K138").Value = 3551.1426
$ws.Range("L138").Value = 10158.96
$ws.Range("M138").Value = 1588.8574
$ws.Range("N138").Value = -20438.96

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2234.8645
$ws.Range("I32").Value = 2120.1553
$ws.Range("K32").Value = 2120.1553
$ws.Range("M32").Value = -1833.1553

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4361.607
$ws.Range("I61").Value = 4338
$ws.Range("J61").Value = 4999
$ws.Range("K61").Value = 4338
$ws.Range("L61").Value = 4999
$ws.Range("M61").Value = -4126
$ws.Range("N61").Value = -5423

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2381.7585
$ws.Range("I74").Value = 2430.75
$ws.Range("J74").Value = 1010
$ws.Range("K74").Value = 2430.75
$ws.Range("L74").Value = 1010
$ws.Range("M74").Value = -1556.75
$ws.Range("N74").Value = -2758

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2381.7585
$ws.Range("I77").Value = 2430.75
$ws.Range("J77").Value = 1010
$ws.Range("K77").Value = 12153.75
$ws.Range("L77").Value = 5050
$ws.Range("M77").Value = -7785.75
$ws.Range("N77").Value = -13786

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3402.1765
$ws.Range("I132").Value = 2451.9167
$ws.Range("J132").Value = 5682.8
$ws.Range("K132").Value = 7355.750100000001
$ws.Range("L132").Value = 17048.4
$ws.Range("M132").Value = -4825.750100000001
$ws.Range("N132").Value = -22108.4

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 77332.664
$ws.Range("J133").Value = 77332.664
$ws.Range("L133").Value = 77332.664
$ws.Range("N133").Value = -82392.664

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4361.607
$ws.Range("I136").Value = 4338
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 13014
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -10464
$ws.Range("N136").Value = -20097

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2978.8696
$ws.Range("I99").Value = 1500.875
$ws.Range("J99").Value = 6357.143
$ws.Range("K99").Value = 1500.875
$ws.Range("L99").Value = 6357.143
$ws.Range("M99").Value = -2.875
$ws.Range("N99").Value = -9353.143

# CRP row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 200479.8
$ws.Range("I4").Value = 200479.8
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 200479.8
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -200367.8
$ws.Range("N4").ClearContents()

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2989.627
$ws.Range("I31").Value = 1588.1562
$ws.Range("J31").Value = 4270.971
$ws.Range("K31").Value = 1588.1562
$ws.Range("L31").Value = 4270.971
$ws.Range("M31").Value = -1293.1562
$ws.Range("N31").Value = -4860.971

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2989.627
$ws.Range("I34").Value = 1588.1562
$ws.Range("J34").Value = 4270.971
$ws.Range("K34").Value = 1588.1562
$ws.Range("L34").Value = 4270.971
$ws.Range("M34").Value = -1386.1562
$ws.Range("N34").Value = -4674.971

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2425.3076
$ws.Range("I58").Value = 1869.8889
$ws.Range("J58").Value = 3675
$ws.Range("K58").Value = 1869.8889
$ws.Range("L58").Value = 3675
$ws.Range("M58").Value = -1666.8889
$ws.Range("N58").Value = -4081

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2222
$ws.Range("I132").Value = 2249.75
$ws.Range("K132").Value = 6749.25
$ws.Range("M132").Value = -4219.25

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2647.1428
$ws.Range("I134").Value = 3024.6
$ws.Range("J134").Value = 1703.5
$ws.Range("K134").Value = 9073.799999999999
$ws.Range("L134").Value = 5110.5
$ws.Range("M134").Value = -6538.799999999999
$ws.Range("N134").Value = -10180.5

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2425.3076
$ws.Range("I136").Value = 1869.8889
$ws.Range("J136").Value = 3675
$ws.Range("K136").Value = 5609.6667
$ws.Range("L136").Value = 11025
$ws.Range("M136").Value = -3059.6667
$ws.Range("N136").Value = -16125

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 60462250
$ws.Range("I4").Value = 30747172
$ws.Range("J4").Value = 226866670
$ws.Range("K4").Value = 92241516
$ws.Range("L4").Value = 680600010
$ws.Range("M4").Value = -92241404
$ws.Range("N4").Value = -680600234

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 803946.6
$ws.Range("I11").Value = 4645.476
$ws.Range("J11").Value = 5000278
$ws.Range("K11").Value = 13936.428
$ws.Range("L11").Value = 15000834
$ws.Range("M11").Value = -13796.428
$ws.Range("N11").Value = -15001114

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4790
$ws.Range("I68").Value = 3450.3333
$ws.Range("K68").Value = 10350.9999
$ws.Range("M68").Value = -9539.999899999999

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 4790
$ws.Range("I71").Value = 3450.3333
$ws.Range("K71").Value = 31052.9997
$ws.Range("M71").Value = -26996.9997

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 942.75
$ws.Range("I113").Value = 725.4761999999999
$ws.Range("J113").Value = 1357.5454
$ws.Range("K113").Value = 2176.4286
$ws.Range("L113").Value = 4072.6362
$ws.Range("M113").Value = -6.428599999999733
$ws.Range("N113").Value = -8412.636200000001

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1795.04
$ws.Range("I102").Value = 1921.619
$ws.Range("J102").Value = 1130.5
$ws.Range("K102").Value = 1921.619
$ws.Range("L102").Value = 1130.5
$ws.Range("M102").Value = -299.6189999999999
$ws.Range("N102").Value = -4374.5

# GSM row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 71666.664
$ws.Range("J134").Value = 71666.664
$ws.Range("L134").Value = 214999.992
$ws.Range("N134").Value = -220069.992

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 120000.336
$ws.Range("J133").Value = 115000.5
$ws.Range("L133").Value = 115000.5
$ws.Range("N133").Value = -120060.5

# LTW row 135
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5419.9
$ws.Range("I136").Value = 4414.143
$ws.Range("J136").Value = 5961.4614
$ws.Range("K136").Value = 13242.429
$ws.Range("L136").Value = 17884.3842
$ws.Range("M136").Value = -10692.429
$ws.Range("N136").Value = -22984.3842

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2961285.8
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2224

# WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 10000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -10226

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2409
$ws.Range("I81").Value = 2022.6
$ws.Range("J81").Value = 3375
$ws.Range("K81").Value = 4045.2
$ws.Range("L81").Value = 6750
$ws.Range("M81").Value = -2984.2
$ws.Range("N81").Value = -8872

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2409
$ws.Range("I84").Value = 2022.6
$ws.Range("J84").Value = 3375
$ws.Range("K84").Value = 20226
$ws.Range("L84").Value = 33750
$ws.Range("M84").Value = -14922
$ws.Range("N84").Value = -44358

# WVR row 105
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 75000
$ws.Range("J105").Value = 75000
$ws.Range("L105").Value = 75000
$ws.Range("N105").Value = -81988

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3104.7778
$ws.Range("I136").Value = 2627.5715
$ws.Range("J136").Value = 4775
$ws.Range("K136").Value = 7882.7145
$ws.Range("L136").Value = 14325
$ws.Range("M136").Value = -5332.7145
$ws.Range("N136").Value = -19425
